# Updated symbol list on Mon Jan 16 14:23:41 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must stay text, matching the
# original inline-string cell type used for this data table).
$updates = @{
    "D2" = "297.14"
    "E2" = "-0.38%"
    "D3" = "31.31"
    "E3" = "-0.44%"
    "D4" = "5.092"
    "E4" = "-1.14%"
    "D5" = "0.07994"
    "E5" = "9.24%"
    "D6" = "2.435"
    "E6" = "30.29%"
    "D7" = "7.788"
    "E7" = "0.34%"
    "D8" = "3.806"
    "E8" = "1.72%"
    "D9" = "0.9223"
    "E9" = "-0.34%"
    "D10" = "0.1729"
    "E10" = "3.15%"
    "D11" = "0.07309"
    "E11" = "2.18%"
    "D12" = "0.09245"
    "E12" = "15.91%"
    "D13" = "0.03039"
    "E13" = "1.39%"
    "D14" = "0.09982"
    "E14" = "0.59%"
    "D15" = "0.001498"
    "E15" = "0.66%"
    "D16" = "0.005905"
    "E16" = "-4.18%"
    "D17" = "3.502"
    "E17" = "1.35%"
    "D18" = "2.246"
    "E18" = "1.22%"
    "D20" = "0.1347"
    "E20" = "2.27%"
    "D21" = "4.589"
    "E21" = "0.84%"
    "D22" = "0.1618"
    "E22" = "2.28%"
    "D23" = "0.04656"
    "E23" = "0.36%"
    "D24" = "0.001243"
    "E24" = "2.22%"
    "D25" = "0.004423"
    "E25" = "-6.60%"
    "E26" = "-7.54%"
    "D27" = "0.0003434"
    "E27" = "83.29%"
    "D39" = "0.01783"
    "E39" = "3.92%"
    "D40" = "0.04446"
    "E40" = "-0.61%"
    "D41" = "0.006965"
    "E41" = "-1.35%"
    "D42" = "0.1341"
    "E42" = "0.81%"
    "D43" = "0.002148"
    "E43" = "-0.47%"
    "D44" = "0.009794"
    "E44" = "-7.45%"
    "D45" = "0.00006613"
    "E45" = "6.14%"
    "D46" = "0.00000000750"
    "E46" = "0.05%"
    "D49" = "0.00002101"
    "E49" = "0.05%"
    "D50" = "0.0002001"
    "E50" = "0.12%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
